# Update "想去人数" (want-to-go count) figures in column F across sheets.
# Sheet "展览" (sheet1): rows 3, 7, 8, 9
# Sheet "演出" (sheet2): row 2
# Sheet "全部类型" (sheet4): rows 3, 7, 8, 9, 11

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 471
$wsExhibit.Range("F7").Value = 36
$wsExhibit.Range("F8").Value = 1098
$wsExhibit.Range("F9").Value = 3878

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 51

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 471
$wsAll.Range("F7").Value = 36
$wsAll.Range("F8").Value = 1098
$wsAll.Range("F9").Value = 3878
$wsAll.Range("F11").Value = 51
